$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.266960666666667
$ws.Range("H2").Value = 6.800882000000001
$ws.Range("I2").Value = 0.0133060725960104
$ws.Range("J2").Value = 0.0133060725960104
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 381.0983662447916
$ws.Range("R2").Value = 3429.885296203124
$ws.Range("S2").Value = 0.003970774576510696
$ws.Range("T2").Value = 0.003970774576510697
$ws.Range("G3").Value = 2.266960666666667
$ws.Range("H3").Value = 6.800882000000001
$ws.Range("I3").Value = 0.0133060725960104
$ws.Range("J3").Value = 0.0133060725960104
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 369.5287246777305
$ws.Range("R3").Value = 3325.758522099574
$ws.Range("S3").Value = 0.003850227120360442
$ws.Range("T3").Value = 0.003850227120360442
$ws.Range("G4").Value = 2.266960666666667
$ws.Range("H4").Value = 6.800882000000001
$ws.Range("I4").Value = 0.0133060725960104
$ws.Range("J4").Value = 0.0133060725960104
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 376.3008034311533
$ws.Range("R4").Value = 3386.70723088038
$ws.Range("S4").Value = 0.003920787375995196
$ws.Range("T4").Value = 0.003920787375995196
$ws.Range("G5").Value = 2.266960666666667
$ws.Range("H5").Value = 6.800882000000001
$ws.Range("I5").Value = 0.0133060725960104
$ws.Range("J5").Value = 0.0133060725960104
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 150.1334018154496
$ws.Range("R5").Value = 1351.200616339046
$ws.Range("S5").Value = 0.001564283523144071
$ws.Range("T5").Value = 0.001564283523144071
$ws.Range("I6").Value = 0.7568089559072322
$ws.Range("J6").Value = 0.7568089559072322
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 21675.7164501079
$ws.Range("R6").Value = 195081.4480509711
$ws.Range("S6").Value = 0.2258455858938478
$ws.Range("T6").Value = 0.2258455858938478
$ws.Range("I7").Value = 0.7568089559072322
$ws.Range("J7").Value = 0.7568089559072322
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.2189892130784987
$ws.Range("T7").Value = 0.2189892130784987
$ws.Range("I8").Value = 0.7568089559072322
$ws.Range("J8").Value = 0.7568089559072322
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 21402.84566290225
$ws.Range("R8").Value = 192625.6109661203
$ws.Range("S8").Value = 0.2230024659004845
$ws.Range("T8").Value = 0.2230024659004845
$ws.Range("I9").Value = 0.7568089559072322
$ws.Range("J9").Value = 0.7568089559072322
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 8539.131457077652
$ws.Range("R9").Value = 76852.18311369888
$ws.Range("S9").Value = 0.08897169103440124
$ws.Range("T9").Value = 0.08897169103440124
$ws.Range("G10").Value = 37.26833833333333
$ws.Range("H10").Value = 111.805015
$ws.Range("I10").Value = 0.2187489278872993
$ws.Range("J10").Value = 0.2187489278872993
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 6265.173922216914
$ws.Range("R10").Value = 56386.56529995223
$ws.Range("S10").Value = 0.06527866695649137
$ws.Range("T10").Value = 0.06527866695649139
$ws.Range("G11").Value = 37.26833833333333
$ws.Range("H11").Value = 111.805015
$ws.Range("I11").Value = 0.2187489278872993
$ws.Range("J11").Value = 0.2187489278872993
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 6074.971541268401
$ws.Range("R11").Value = 54674.7438714156
$ws.Range("S11").Value = 0.0632968931008222
$ws.Range("T11").Value = 0.0632968931008222
$ws.Range("G12").Value = 37.26833833333333
$ws.Range("H12").Value = 111.805015
$ws.Range("I12").Value = 0.2187489278872993
$ws.Range("J12").Value = 0.2187489278872993
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 6186.303037184316
$ws.Range("R12").Value = 55676.72733465885
$ws.Range("S12").Value = 0.0644568882955113
$ws.Range("T12").Value = 0.0644568882955113
$ws.Range("G13").Value = 37.26833833333333
$ws.Range("H13").Value = 111.805015
$ws.Range("I13").Value = 0.2187489278872993
$ws.Range("J13").Value = 0.2187489278872993
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 2468.16034184645
$ws.Range("R13").Value = 22213.44307661804
$ws.Range("S13").Value = 0.02571647953447445
$ws.Range("T13").Value = 0.02571647953447445
$ws.Range("G14").Value = 1.897252
$ws.Range("H14").Value = 5.691756
$ws.Range("I14").Value = 0.01113604360945798
$ws.Range("J14").Value = 0.01113604360945798
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 318.9467061278213
$ws.Range("R14").Value = 2870.520355150392
$ws.Range("S14").Value = 0.003323198376402092
$ws.Range("T14").Value = 0.003323198376402092
$ws.Range("G15").Value = 1.897252
$ws.Range("H15").Value = 5.691756
$ws.Range("I15").Value = 0.01113604360945798
$ws.Range("J15").Value = 0.01113604360945798
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 309.2639066310547
$ws.Range("R15").Value = 2783.375159679492
$ws.Range("S15").Value = 0.003222310475858023
$ws.Range("T15").Value = 0.003222310475858023
$ws.Range("G16").Value = 1.897252
$ws.Range("H16").Value = 5.691756
$ws.Range("I16").Value = 0.01113604360945798
$ws.Range("J16").Value = 0.01113604360945798
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 314.93155677956
$ws.Range("R16").Value = 2834.38401101604
$ws.Range("S16").Value = 0.003281363369051972
$ws.Range("T16").Value = 0.003281363369051972
$ws.Range("G17").Value = 1.897252
$ws.Range("H17").Value = 5.691756
$ws.Range("I17").Value = 0.01113604360945798
$ws.Range("J17").Value = 0.01113604360945798
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 125.6488041673853
$ws.Range("R17").Value = 1130.839237506468
$ws.Range("S17").Value = 0.001309171388145891
$ws.Range("T17").Value = 0.001309171388145891
